$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.490.99'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '3.114.52'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("D5").Value = '''526.74'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").Value = '''137.41'
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("D8").Value = '3.115.74'
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D9").Value = '''0.447'
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("D12").Value = '''0.397'
$ws.Range("E12").Value = '  +3.22%  '
$ws.Range("D13").Value = '3.652.11'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '57.637.26'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = '3.114.09'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").Value = '''5.96'
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").Value = '''12.61'
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("E21").Value = '  -1.78%  '
$ws.Range("D22").Value = '''348.75'
$ws.Range("E22").Value = '  +3.78%  '
$ws.Range("D23").Value = '''5.79'
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '''68.35'
$ws.Range("E25").Value = '  +2.67%  '
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '0.0₃0914'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").Value = '''7.43'
$ws.Range("E30").Value = '  +3.29%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").Value = '''6.10'
$ws.Range("E33").Value = '  -6.46%  '
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("E36").Value = '  +6.98%  '
$ws.Range("D37").Value = '''158.43'
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("E38").Value = '  +1.13%  '
$ws.Range("D39").Value = '''26.14'
$ws.Range("E39").Value = '  -3.06%  '
$ws.Range("E40").Value = '  -2.52%  '
$ws.Range("E41").Value = '  +6.81%  '
$ws.Range("E42").Value = '  +0.64%  '
$ws.Range("E43").Value = '  +7.24%  '
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").Value = '3.156.10'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").Value = '''36.47'
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.345.61'
$ws.Range("E47").Value = '  +1.93%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '''0.0269'
$ws.Range("E48").Value = '  +3.66%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = '''0.960'
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("E51").Value = '  +0.46%  '

# Strip the quote-prefix style artifact introduced by forcing text above
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D50").ClearFormats()
